# Update column F ("dSF") values for several rows, as per re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 2
$ws.Range("F4").Value = -6
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = 1
$ws.Range("F10").Value = -4
$ws.Range("F11").Value = 2
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = 3
$ws.Range("F16").Value = 1
